# Modify Runaway condition and regenerated health
# Inserts a new "sprintSpeed" column before the "stoppingDist" column (old column N),
# fills in sprintSpeed values for Hunter / Bot_X / Player, and bumps the
# "acceleration" column (now P) from 2 to 8 for those same rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PropertyInfo")

# Insert a new column at N; everything from N..R shifts right to O..S
$ws.Columns("N").Insert()

# New column header (row 2) for the inserted "sprintSpeed" stat
$ws.Range("N2").Value = "sprintSpeed"

# New column type marker (row 3) matches the other "!Float" columns
$ws.Range("N3").Value = "!Float"

# New sprintSpeed values for the data rows (Hunter, Bot_X, Player)
$ws.Range("N6").Value = 5.1
$ws.Range("N7").Value = 4
$ws.Range("N8").Value = 3.5

# acceleration column (now P, was O before insert) values bumped from 2 to 8
$ws.Range("P6").Value = 8
$ws.Range("P7").Value = 8
$ws.Range("P8").Value = 8

# Selection moved to N7 and the top-left frozen cell reset (no longer E1)
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("N7").Select()
